{"js": "const replacements = [\n  [\"[Music]\", \"[Musica]\"],\n  [\"Introduction\", \"Introduzione\"],\n  [\"Hi everyone I'm sony from London in the\", \"Ciao a tutti, sono Sony, da Londra, Regno\"],\n  [\"UK and today I have a really exciting\", \"Unito e oggi ho un video davvero\"],\n  [\"video for you for the Virtual Maths\", \"entusiasmante per voi per il Campo Virtuale\"],\n  [\"Camp (VMC)\", \"di Matematica (VMC)\"],\n  [\"We're going to play a fun mathematical game called split or steal\", \"Giocheremo a un divertente gioco matematico noto come dividi o ruba,\"],\n  [\"and this game is a two player game.\", \"un gioco per due giocatori.\"],\n  [\"So a quick introduction split or steal\", \"Quindi, ecco una rapida introduzione a dividi o ruba\"],\n  [\"is a game based on the famous Prisoners Dilemma which you might have\", \"\u00c8 un gioco basato sul famoso Dilemma dei Prigionieri, che potreste\"],\n  [\"heard of before\", \"aver sentito prima\"],\n  [\"It's an interesting new branch of maths called Game Theory\", \"\u00c8 un interessante nuovo ramo della matematica, detto Teoria del Gioco,\"],\n  [\"Where when making your choice you also need to consider\", \"in cui facendo le proprie scelte, si deve anche considerare\"],\n  [\"The other player's choice\", \"la scelta dell'altro giocatore\"],\n  [\"This means your choices are interdependent \", \"Questo significa che le tue scelte sono interdipendenti \"],\n  [\"So what you need today: a partner to play the game with,\", \"Quindi, ecco cosa ti serve oggi: un partner con cui giocare,\"],\n  [\"two small cards each saying split and steal,\", \"due piccole carte che dicano ognuna dividi e ruba,\"],\n  [\"and these can just be small bits of paper\", \"che possono anche essere piccoli pezzi di carta\"],\n  [\"and finally something to count scores with\", \"e, infine, qualcosa con cui contare i punteggi\"],\n  [\"If you have some tokens or coins or something like that,\", \"Se hai dei gettoni o delle monete o qualcosa di simile,\"],\n  [\"That would be great,\", \"sarebbe fantastico,\"],\n  [\"but if not you could just use a pen and paper to write down the\", \"altrimenti, basta usare carta e penna per annotare i\"],\n  [\"scores \", \"punteggi \"],\n  [\"So this is what the game looks and this table is called a matrix\", \"Quindi, ecco come appare il gioco e, questa tabella, si chiama matrice\"],\n  [\"In the game there are two points to be won but who wins what is decided by the choices of the players\", \"Nel gioco ci sono due punti da vincere, ma chi vince cosa \u00e8 deciso dalle scelte dei giocatori\"],\n  [\"We have our two players on the left and at the top\", \"Abbiamo i nostri due giocatori sulla sinistra e in alto:\"],\n  [\"The red player and the blue player\", \"il giocatore rosso e il giocatore blu;\"],\n  [\"each player has two options shown next to them, split or steal\", \"ogni giocatore ha due opzioni, una affianco all'altra, dividere o rubare\"],\n  [\"Since two players have two choices each there are four outcomes in total\", \"Poich\u00e9 i due giocatori hanno due scelte l'uno, ci sono quattro risultati in tutto\"],\n  [\"and they are all shown in the table in each section\", \"e sono tutti mostrati nella tabella in ogni sezione\"],\n  [\"The red number is the number of points won by the red player\", \"Il numero rosso \u00e8 il numero di punti vinto dal giocatore rosso,\"],\n  [\"and the blue number is the number of points won by the blue player\", \"il numero blu \u00e8 il numero di punti vinto dal giocatore blu\"],\n  [\"for example\", \"ad esempio\"],\n  [\"if both players choose to split\", \"se entrambi i giocatori decidono di dividere\"],\n  [\"we would end up with the top left outcome\", \"finiremmo con il risultato in alto a sinistra\"],\n  [\"and the players would split the two points to earn one point each\", \"e i giocatori dividerebbero i due punti per ottenere un punto ciascuno.\"],\n  [\"However if the red player wanted to split\", \"Tuttavia, se il giocatore rosso avesse voluto dividere\"],\n  [\"But the blue player chose to steal the blue player would steal the red player's point and earn two points\", \"ma il giocatore blu avesse scelto di rubare, il giocatore blu avrebbe rubato il punto del giocatore rosso e vinto due punti;\"],\n  [\"while the red player wins nothing\", \"mentre il giocatore rosso non vince nulla\"],\n  [\"The opposite happens if the blue player splits and the red player steals\", \"L'opposto si verifica se il giocatore blu divide e quello rosso ruba,\"],\n  [\"but if both players try and steal\", \"ma se entrambi provano a rubare\"],\n  [\"it doesn't work and no one wins the points so both players end up with nothing\", \"non funziona, e nessuno vince i punti, quindi entrambi i giocatori finiscono senza niente.\"],\n  [\"Now it's your turn get ready to play split or steal with your partner\", \"Ora tocca a te, preparati a giocare a dividi o ruba con il tuo partner.\"],\n  [\"First talk to your partner for a couple of minutes about what choice you're going to make\", \"Prima, parla con il tuo partner per qualche minuto su quale scelta effettuerai.\"],\n  [\"Remember, the person with the most points wins\", \"Ricorda, la persona con pi\u00f9 punti vince\"],\n  [\"You are allowed to lie to your partner then secretly choose split or steal\", \"Puoi mentire al tuo partner e poi scegliere di dividere o rubare in segreto\"],\n  [\"and place the card you have chosen face down so your partner can't see it\", \"e posiziona la carta che hai scelto a faccia in gi\u00f9, cos\u00ec che il tuo partner non possa vederla.\"],\n  [\"Finally reveal your choices and work out your scores\", \"Infine, rivela le tue scelte e calcola i tuoi punteggi\"],\n  [\"Play the game once with your partner\", \"Gioca una volta con il tuo partner;\"],\n  [\"Pause the video now\", \"ora, interrompi il video.\"],\n  [\"How did it go? \", \"Com'\u00e8 andata? \"],\n  [\"Did you get the number of points you were hoping for?\", \"Hai ricevuto il numero di punti che speravi?\"],\n  [\"did you and your partner tell the truth to each other?\", \"Tu e il tuo partner vi siete detti la verit\u00e0?\"],\n  [\"Let's think about why the result might have been different to what you expected\", \"Pensiamo al perch\u00e9 il risultato potrebbe esser stato differente da quanto ti aspettavi\"],\n  [\"Imagine you are the red player your opponent\", \"Immagina di essere il giocatore rosso, tuo avversario.\"],\n  [\"The blue player has two choices split or steal \", \"Il giocatore blu ha due scelte: dividere o rubare; \"],\n  [\"if blue chooses split\", \"se il blu sceglie di dividere\"],\n  [\"you could either choose split and win one point\", \"potresti scegliere di dividere e vincere un punto\"],\n  [\"or you could choose steal and win two points\", \"o di rubare e vincerne due.\"],\n  [\"two points is better than one so you would choose steal\", \"Due punti sono meglio di uno quindi sceglieresti di rubare\"],\n  [\"What if the blue player chose steal?\", \"E se il giocatore blu avesse scelto di rubare?\"],\n  [\"If you choose split you get zero and if you choose steal you also get zero\", \"Se scegli di dividere ottieni zero e se scegli di rubare, ottieni comunque zero.\"],\n  [\"So it doesn't really matter what you choose\", \"Quindi non importa davvero cosa scegli\"],\n  [\"But let's assume you prefer to steal so you don't give your opponent any points\", \"Ma presumiamo che preferisci rubare, cos\u00ec da non dare alcun punto al tuo avversario.\"],\n  [\"As we have seen no matter what your opponent does\", \"Come abbiamo visto, non importa cosa il tuo avversario faccia,\"],\n  [\"split is never the best choice this means steal is called a weakly dominant strategy\", \"dividere non \u00e8 mai la scelta migliore; questo significa che rubare \u00e8 una strategia debolmente dominante.\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n$replacements = @(\n    ,@(\"[Music]\", \"[Musica]\")\n    ,@(\"Introduction\", \"Introduzione\")\n    ,@(\"Hi everyone I'm sony from London in the\", \"Ciao a tutti, sono Sony, da Londra, Regno\")\n    ,@(\"UK and today I have a really exciting\", \"Unito e oggi ho un video davvero\")\n    ,@(\"video for you for the Virtual Maths\", \"entusiasmante per voi per il Campo Virtuale\")\n    ,@(\"Camp (VMC)\", \"di Matematica (VMC)\")\n    ,@(\"We're going to play a fun mathematical game called split or steal\", \"Giocheremo a un divertente gioco matematico noto come dividi o ruba,\")\n    ,@(\"and this game is a two player game.\", \"un gioco per due giocatori.\")\n    ,@(\"So a quick introduction split or steal\", \"Quindi, ecco una rapida introduzione a dividi o ruba\")\n    ,@(\"is a game based on the famous Prisoners Dilemma which you might have\", \"\u00c8 un gioco basato sul famoso Dilemma dei Prigionieri, che potreste\")\n    ,@(\"heard of before\", \"aver sentito prima\")\n    ,@(\"It's an interesting new branch of maths called Game Theory\", \"\u00c8 un interessante nuovo ramo della matematica, detto Teoria del Gioco,\")\n    ,@(\"Where when making your choice you also need to consider\", \"in cui facendo le proprie scelte, si deve anche considerare\")\n    ,@(\"The other player's choice\", \"la scelta dell'altro giocatore\")\n    ,@(\"This means your choices are interdependent \", \"Questo significa che le tue scelte sono interdipendenti \")\n    ,@(\"So what you need today: a partner to play the game with,\", \"Quindi, ecco cosa ti serve oggi: un partner con cui giocare,\")\n    ,@(\"two small cards each saying split and steal,\", \"due piccole carte che dicano ognuna dividi e ruba,\")\n    ,@(\"and these can just be small bits of paper\", \"che possono anche essere piccoli pezzi di carta\")\n    ,@(\"and finally something to count scores with\", \"e, infine, qualcosa con cui contare i punteggi\")\n    ,@(\"If you have some tokens or coins or something like that,\", \"Se hai dei gettoni o delle monete o qualcosa di simile,\")\n    ,@(\"That would be great,\", \"sarebbe fantastico,\")\n    ,@(\"but if not you could just use a pen and paper to write down the\", \"altrimenti, basta usare carta e penna per annotare i\")\n    ,@(\"scores \", \"punteggi \")\n    ,@(\"So this is what the game looks and this table is called a matrix\", \"Quindi, ecco come appare il gioco e, questa tabella, si chiama matrice\")\n    ,@(\"In the game there are two points to be won but who wins what is decided by the choices of the players\", \"Nel gioco ci sono due punti da vincere, ma chi vince cosa \u00e8 deciso dalle scelte dei giocatori\")\n    ,@(\"We have our two players on the left and at the top\", \"Abbiamo i nostri due giocatori sulla sinistra e in alto:\")\n    ,@(\"The red player and the blue player\", \"il giocatore rosso e il giocatore blu;\")\n    ,@(\"each player has two options shown next to them, split or steal\", \"ogni giocatore ha due opzioni, una affianco all'altra, dividere o rubare\")\n    ,@(\"Since two players have two choices each there are four outcomes in total\", \"Poich\u00e9 i due giocatori hanno due scelte l'uno, ci sono quattro risultati in tutto\")\n    ,@(\"and they are all shown in the table in each section\", \"e sono tutti mostrati nella tabella in ogni sezione\")\n    ,@(\"The red number is the number of points won by the red player\", \"Il numero rosso \u00e8 il numero di punti vinto dal giocatore rosso,\")\n    ,@(\"and the blue number is the number of points won by the blue player\", \"il numero blu \u00e8 il numero di punti vinto dal giocatore blu\")\n    ,@(\"for example\", \"ad esempio\")\n    ,@(\"if both players choose to split\", \"se entrambi i giocatori decidono di dividere\")\n    ,@(\"we would end up with the top left outcome\", \"finiremmo con il risultato in alto a sinistra\")\n    ,@(\"and the players would split the two points to earn one point each\", \"e i giocatori dividerebbero i due punti per ottenere un punto ciascuno.\")\n    ,@(\"However if the red player wanted to split\", \"Tuttavia, se il giocatore rosso avesse voluto dividere\")\n    ,@(\"But the blue player chose to steal the blue player would steal the red player's point and earn two points\", \"ma il giocatore blu avesse scelto di rubare, il giocatore blu avrebbe rubato il punto del giocatore rosso e vinto due punti;\")\n    ,@(\"while the red player wins nothing\", \"mentre il giocatore rosso non vince nulla\")\n    ,@(\"The opposite happens if the blue player splits and the red player steals\", \"L'opposto si verifica se il giocatore blu divide e quello rosso ruba,\")\n    ,@(\"but if both players try and steal\", \"ma se entrambi provano a rubare\")\n    ,@(\"it doesn't work and no one wins the points so both players end up with nothing\", \"non funziona, e nessuno vince i punti, quindi entrambi i giocatori finiscono senza niente.\")\n    ,@(\"Now it's your turn get ready to play split or steal with your partner\", \"Ora tocca a te, preparati a giocare a dividi o ruba con il tuo partner.\")\n    ,@(\"First talk to your partner for a couple of minutes about what choice you're going to make\", \"Prima, parla con il tuo partner per qualche minuto su quale scelta effettuerai.\")\n    ,@(\"Remember, the person with the most points wins\", \"Ricorda, la persona con pi\u00f9 punti vince\")\n    ,@(\"You are allowed to lie to your partner then secretly choose split or steal\", \"Puoi mentire al tuo partner e poi scegliere di dividere o rubare in segreto\")\n    ,@(\"and place the card you have chosen face down so your partner can't see it\", \"e posiziona la carta che hai scelto a faccia in gi\u00f9, cos\u00ec che il tuo partner non possa vederla.\")\n    ,@(\"Finally reveal your choices and work out your scores\", \"Infine, rivela le tue scelte e calcola i tuoi punteggi\")\n    ,@(\"Play the game once with your partner\", \"Gioca una volta con il tuo partner;\")\n    ,@(\"Pause the video now\", \"ora, interrompi il video.\")\n    ,@(\"How did it go? \", \"Com'\u00e8 andata? \")\n    ,@(\"Did you get the number of points you were hoping for?\", \"Hai ricevuto il numero di punti che speravi?\")\n    ,@(\"did you and your partner tell the truth to each other?\", \"Tu e il tuo partner vi siete detti la verit\u00e0?\")\n    ,@(\"Let's think about why the result might have been different to what you expected\", \"Pensiamo al perch\u00e9 il risultato potrebbe esser stato differente da quanto ti aspettavi\")\n    ,@(\"Imagine you are the red player your opponent\", \"Immagina di essere il giocatore rosso, tuo avversario.\")\n    ,@(\"The blue player has two choices split or steal \", \"Il giocatore blu ha due scelte: dividere o rubare; \")\n    ,@(\"if blue chooses split\", \"se il blu sceglie di dividere\")\n    ,@(\"you could either choose split and win one point\", \"potresti scegliere di dividere e vincere un punto\")\n    ,@(\"or you could choose steal and win two points\", \"o di rubare e vincerne due.\")\n    ,@(\"two points is better than one so you would choose steal\", \"Due punti sono meglio di uno quindi sceglieresti di rubare\")\n    ,@(\"What if the blue player chose steal?\", \"E se il giocatore blu avesse scelto di rubare?\")\n    ,@(\"If you choose split you get zero and if you choose steal you also get zero\", \"Se scegli di dividere ottieni zero e se scegli di rubare, ottieni comunque zero.\")\n    ,@(\"So it doesn't really matter what you choose\", \"Quindi non importa davvero cosa scegli\")\n    ,@(\"But let's assume you prefer to steal so you don't give your opponent any points\", \"Ma presumiamo che preferisci rubare, cos\u00ec da non dare alcun punto al tuo avversario.\")\n    ,@(\"As we have seen no matter what your opponent does\", \"Come abbiamo visto, non importa cosa il tuo avversario faccia,\")\n    ,@(\"split is never the best choice this means steal is called a weakly dominant strategy\", \"dividere non \u00e8 mai la scelta migliore; questo significa che rubare \u00e8 una strategia debolmente dominante.\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $range.Find.MatchCase = $true\n    $range.Find.MatchWholeWord = $false\n    $range.Find.MatchWildcards = $false\n    $range.Find.Execute($old)\n    if ($range.Find.Found) {\n        $range.Text = $new\n    }\n}\n"}
